$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (Through 2022-03-24 -> Through 2022-03-25)
$ws.Name = "Through 2022-03-25"

# Update the header label in I1 (shared string "2022 (through 03-24)" -> "2022 (through 03-25)")
$ws.Range("I1").Value = "2022 (through 03-25)"

# Update the March value (I4): 99 -> 102
$ws.Range("I4").Value = 102

# Update the Total value (I14): 399 -> 402
$ws.Range("I14").Value = 402
